# Actualizo EMAE Junio24 e ICA julio 24
# Updates EMAE (June 2024 row added), Expo-ICA, Impo-ICA and BC por zonas data,
# and moves the active/selected tab from "Aperturas" to "BC por zonas".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) EMAE sheet: revise existing monthly index values (rows 2-102) and append
#    the new June-2024 row (row 103).
# ---------------------------------------------------------------------------
$emaeData = @(
    @{Row=2; B=$null; C=147.97285936103225; D=147.08983858006528},
    @{Row=3; B=$null; C=146.97171661594652; D=146.52101725643882},
    @{Row=4; B=$null; C=146.19934307368629; D=145.97996952659415},
    @{Row=5; B=$null; C=144.99327796692558; D=145.50097797460367},
    @{Row=6; B=$null; C=144.49596830437281; D=145.11066141097464},
    @{Row=7; B=$null; C=144.18660236957894; D=144.82913372993175},
    @{Row=8; B=$null; C=144.41779036980932; D=144.66858777921993},
    @{Row=9; B=$null; C=145.520639999829; D=144.63666451266212},
    @{Row=10; B=$null; C=144.91190335582201; D=144.73751697014367},
    @{Row=11; B=$null; C=145.00687443204001; D=144.96680805047811},
    @{Row=12; B=$null; C=145.82428771823814; D=145.31253413670319},
    @{Row=13; B=$null; C=147.07008993784754; D=145.76221520983762},
    @{Row=14; B=$null; C=147.23816076572268; D=146.2992067382566},
    @{Row=15; B=132.15851633982271; C=146.81664079258277; D=146.90273954702866},
    @{Row=16; B=$null; C=147.80662258063549; D=147.55060877670735},
    @{Row=17; B=$null; C=147.62397867673587; D=148.21856197040509},
    @{Row=18; B=168.38920946936776; C=148.46833957307956; D=148.87728834843458},
    @{Row=19; B=$null; C=150.1213125212825; D=149.49334726861744},
    @{Row=20; B=150.30605012391993; C=150.40230626118557; D=150.03082005025666},
    @{Row=21; B=$null; C=150.45168072138901; D=150.45365983489407},
    @{Row=22; B=146.38655965775379; C=151.3045466620421; D=150.72858699188299},
    @{Row=23; B=149.38594966601448; C=151.80740181380838; D=150.82975470357661},
    @{Row=24; B=$null; C=152.56539482479721; D=150.73550182189106},
    @{Row=25; B=146.78338490922448; C=152.22031886285021; D=150.43747705043666},
    @{Row=26; B=142.74091260617226; C=150.79608454340178; D=149.94500725344884},
    @{Row=27; B=138.8180403516574; C=151.72665365395162; D=149.28184129303918},
    @{Row=28; B=$null; C=151.2055358722649; D=148.48444951665456},
    @{Row=29; B=$null; C=146.80276429211537; D=147.59697323491261},
    @{Row=30; B=$null; C=144.70817680650003; D=146.67117834450858},
    @{Row=31; B=$null; C=143.53403141198291; D=145.76097025974022},
    @{Row=32; B=145.96352443000646; C=143.89225993950862; D=144.91577713861136},
    @{Row=33; B=$null; C=146.68276932286082; D=144.1801838976323},
    @{Row=34; B=137.74656971864525; C=143.23224151075175; D=143.58505527520219},
    @{Row=35; B=142.84327598455894; C=143.63904333406268; D=143.14455672893359},
    @{Row=36; B=140.59240732533195; C=141.6830805367992; D=142.86027157687724},
    @{Row=37; B=136.25161596906852; C=141.89398427580869; D=142.72158508029278},
    @{Row=38; B=$null; C=141.83593525124556; D=142.70117971643376},
    @{Row=39; B=132.26788861275296; C=143.9511013579442; D=142.76158934061928},
    @{Row=40; B=144.96325495569374; C=142.34097241138019; D=142.86162573770625},
    @{Row=41; B=149.91622140334627; C=142.62750151806264; D=142.95443335343393},
    @{Row=42; B=164.13569907584596; C=144.57968762508233; D=142.99868652248597},
    @{Row=43; B=150.85897174138188; C=143.82199612623558; D=142.95973015703774},
    @{Row=44; B=146.77702964086899; C=145.62061119403404; D=142.80895986679903},
    @{Row=45; B=$null; C=144.62644904408546; D=142.53029051596641},
    @{Row=46; B=134.877066478801; C=140.42503738736158; D=142.12017456021849},
    @{Row=47; B=141.63933661339601; C=143.84685277719436; D=141.58670152417045},
    @{Row=48; B=137.7718296678064; C=141.01507399685477; D=140.94847607142862},
    @{Row=49; B=135.76515453277773; C=140.09440852145477; D=140.23646720434061},
    @{Row=50; B=133.89108610595787; C=140.59324354125502; D=139.48581829194953},
    @{Row=51; B=128.97363875350521; C=138.90014775179918; D=138.73492873394517},
    @{Row=52; B=128.211060716964; C=125.54697807968385; D=138.02394302675083},
    @{Row=53; B=113.29503437158793; C=106.20903097722967; D=137.38837726577577},
    @{Row=54; B=131.02956916606541; C=117.65166311769086; D=136.8587831461993},
    @{Row=55; B=132.52196537198478; C=124.86906992288233; D=136.45757538253014},
    @{Row=56; B=127.4636898059227; C=126.26965785540263; D=136.1989640027287},
    @{Row=57; B=125.18389390153655; C=128.91873015516268; D=136.09217700972508},
    @{Row=58; B=127.17507561553394; C=130.79557777810396; D=136.13716942776668},
    @{Row=59; B=131.34550037212989; C=133.44074098855805; D=136.32930055682539},
    @{Row=60; B=132.06865461125233; C=134.85681238005697; D=136.65688286501779},
    @{Row=61; B=133.85436775874837; C=136.9618839586459; D=137.10527156268313},
    @{Row=62; B=131.52154429258059; C=139.58509672538554; D=137.65450365716865},
    @{Row=63; B=126.23926987495938; C=137.49709239445235; D=138.28777379502989},
    @{Row=64; B=145.95689030179008; C=140.55355272209928; D=138.9891681872823},
    @{Row=65; B=147.28081278866429; C=139.5936606514135; D=139.7467484462334},
    @{Row=66; B=151.16930917416479; C=139.0274002789111; D=140.55165186670456},
    @{Row=67; B=148.97966150309085; C=141.71244924700468; D=141.39603276151701},
    @{Row=68; B=142.42605450804857; C=141.54428555226775; D=142.27580045906214},
    @{Row=69; B=140.97500035295286; C=143.49523272225446; D=143.18927585275426},
    @{Row=70; B=141.29588801741792; C=144.06024194023152; D=144.12940603631904},
    @{Row=71; B=139.51411805324415; C=143.36912315951815; D=145.08349126097238},
    @{Row=72; B=143.75160756388769; C=145.86292561784524; D=146.03089495344574},
    @{Row=73; B=147.23078901956941; C=150.0398843895662; D=146.94354820668113},
    @{Row=74; B=139.46393626358665; C=147.95800451137202; D=147.79001568885857},
    @{Row=75; B=138.00627026935322; C=150.12507033152397; D=148.53835719638252},
    @{Row=76; B=153.92755964196547; C=149.53169105041306; D=149.15910063551772},
    @{Row=77; B=156.07720261283845; C=150.62333171105931; D=149.62773055625374},
    @{Row=78; B=163.11692077399596; C=150.65378818012496; D=149.93102231035826},
    @{Row=79; B=159.60238517356794; C=152.20483396636769; D=150.06707510010065},
    @{Row=80; B=151.36804047012006; C=151.60275060258959; D=150.0483596925163},
    @{Row=81; B=150.48692368931995; C=151.02911126438323; D=149.90002688169463},
    @{Row=82; B=148.3062125789998; C=149.85130703722137; D=149.65526040443925},
    @{Row=83; B=144.90674880309797; C=147.8347794918019; D=149.35139364511986},
    @{Row=84; B=146.50322679719753; C=147.25814703857731; D=149.02257922247225},
    @{Row=85; B=144.49763294297557; C=147.59024477954043; D=148.69656450882485},
    @{Row=86; B=143.02671289635725; C=149.503309249799; D=148.39167727353291},
    @{Row=87; B=137.6033896154465; C=149.43962201171203; D=148.11843253907395},
    @{Row=88; B=155.36186621144915; C=150.34561896379077; D=147.87093820398442},
    @{Row=89; B=149.13923794434029; C=146.689318212042; D=147.63327516185882},
    @{Row=90; B=152.72484714790204; C=144.54236448504074; D=147.38222646553785},
    @{Row=91; B=151.64704872678965; C=144.95161424274042; D=147.09576126672712},
    @{Row=92; B=148.95277389848641; C=147.86201062584254; D=146.75296823480446},
    @{Row=93; B=150.7284255960451; C=149.39921508011037; D=146.33808630604784},
    @{Row=94; B=147.54696378142793; C=148.68388138025381; D=145.8441146047665},
    @{Row=95; B=146.67440353083421; C=147.59001514798695; D=145.27400815697089},
    @{Row=96; B=145.52851992116507; C=145.55526300739729; D=144.63824663873672},
    @{Row=97; B=138.39104375166951; C=142.76300057947222; D=143.95933440720702},
    @{Row=98; B=137.61813258937534; C=143.13839064886025; D=143.26397351372069},
    @{Row=99; B=133.88349825269793; C=142.82533557643177; D=142.57900206042444},
    @{Row=100; B=142.27716473080739; C=140.6482299778977; D=141.92631744988179},
    @{Row=101; B=146.01389450203277; C=139.36876392197809; D=141.32295894787674},
    @{Row=102; B=155.57318339451891; C=140.30269583553874; D=140.77682610011465},
    @{Row=103; B=145.69606358282977; C=139.81587401285503; D=140.28919739279141}
)

$emaeSheet = $wb.Worksheets.Item("EMAE")
$emaeSheet.Activate()

# Add the new row by copying the formatting (styles/borders) of the last
# existing data row (102) down into the new row (103).
$emaeSheet.Range("A102:D102").Copy()
$emaeSheet.Range("A103:D103").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Date for the new June-2024 row.
$emaeSheet.Range("A103").Value = 45444

foreach ($d in $emaeData) {
    if ($null -ne $d.B) { $emaeSheet.Cells.Item($d.Row, 2).Value = $d.B }
    if ($null -ne $d.C) { $emaeSheet.Cells.Item($d.Row, 3).Value = $d.C }
    if ($null -ne $d.D) { $emaeSheet.Cells.Item($d.Row, 4).Value = $d.D }
}

# Selection moves to the freshly-added row.
$emaeSheet.Range("A103:D103").Select()

# ---------------------------------------------------------------------------
# 2) BC por zonas: revise source columns B & C (column D is the B-C formula
#    and recalculates automatically).
# ---------------------------------------------------------------------------
$bcZonasData = @(
    @{Row=2; B=9208.7319548399992; C=9993.4383400400002},
    @{Row=3; B=7312.46437504; C=7338.6902749600004},
    @{Row=4; B=6330.1850608499999; C=1323.14198296},
    @{Row=5; B=3486.8069678000002; C=358.43163773999999},
    @{Row=6; B=1429.8258285500001; C=117.64612929},
    @{Row=7; B=4708.0592301999995; C=4924.5388798200001},
    @{Row=8; B=3565.7635654300002; C=3849.5062278400001},
    @{Row=9; B=4359.1029573300002; C=5156.6580138899999},
    @{Row=10; B=4067.92805021; C=5608.0038547300001},
    @{Row=11; B=1878.0212924099999; C=736.60161417999996},
    @{Row=12; B=4119.9789890700004; C=1610.3699549999999},
    @{Row=13; B=1887.8990394800001; C=321.79855132},
    @{Row=14; B=897.88899185000002; C=144.60225226},
    @{Row=15; B=2566.1273171500002; C=366.28708497000002},
    @{Row=16; B=1499.33051525; C=550.78583059000005},
    @{Row=17; B=185.31417325999999; C=94.548970280000006},
    @{Row=18; B=340.66265263999998; C=198.00187292000001},
    @{Row=19; B=6133.7311456500001; C=2572.8521916099999}
)

$bcSheet = $wb.Worksheets.Item("BC por zonas")
foreach ($d in $bcZonasData) {
    $bcSheet.Cells.Item($d.Row, 2).Value = $d.B
    $bcSheet.Cells.Item($d.Row, 3).Value = $d.C
}

# ---------------------------------------------------------------------------
# 3) Expo-ICA: revise column B values (ICA julio 24 update).
# ---------------------------------------------------------------------------
$expoIcaData = @(
    @{Row=2; B=45397.173338859997},
    @{Row=3; B=11822.688913739999},
    @{Row=4; B=17.680916029999999},
    @{Row=5; B=999.28606648000004},
    @{Row=6; B=106.77112605000001},
    @{Row=7; B=364.65475953999999},
    @{Row=8; B=396.68078197},
    @{Row=9; B=7208.6130799800003},
    @{Row=10; B=2112.9205309499998},
    @{Row=11; B=150.76669293},
    @{Row=12; B=15.202928679999999},
    @{Row=13; B=83.972491469999994},
    @{Row=14; B=242.79223877000001},
    @{Row=15; B=123.34730089},
    @{Row=16; B=16509.768215169999},
    @{Row=17; B=1933.1152773599999},
    @{Row=18; B=156.78683337999999},
    @{Row=19; B=681.18236092999996},
    @{Row=20; B=29.26410409},
    @{Row=21; B=66.886377780000004},
    @{Row=22; B=108.18168507999999},
    @{Row=23; B=514.62024569000005},
    @{Row=24; B=4297.8076938699996},
    @{Row=25; B=191.48685047999999},
    @{Row=26; B=545.19434192999995},
    @{Row=27; B=455.80783958000001},
    @{Row=28; B=6777.3628022900002},
    @{Row=29; B=123.61610478999999},
    @{Row=30; B=220.7404904},
    @{Row=31; B=55.730753079999999},
    @{Row=32; B=351.98445443999998},
    @{Row=33; B=11400.170645779999},
    @{Row=34; B=2385.1396657199998},
    @{Row=35; B=508.99806318999998},
    @{Row=36; B=123.74743546000001},
    @{Row=37; B=11.178694269999999},
    @{Row=38; B=222.68245612000001},
    @{Row=39; B=57.675819429999997},
    @{Row=40; B=4.3755246799999998},
    @{Row=41; B=78.09677834},
    @{Row=42; B=1510.5992219899999},
    @{Row=43; B=1316.48493717},
    @{Row=44; B=760.82913717999998},
    @{Row=45; B=4208.0772998900002},
    @{Row=46; B=14.51846705},
    @{Row=47; B=197.76714529},
    @{Row=48; B=5664.54556417},
    @{Row=49; B=3114.6785103900002},
    @{Row=50; B=1522.56420065},
    @{Row=51; B=40.91301052},
    @{Row=52; B=843.42786185},
    @{Row=53; B=142.96198075999999}
)

$expoSheet = $wb.Worksheets.Item("Expo-ICA")
foreach ($d in $expoIcaData) {
    $expoSheet.Cells.Item($d.Row, 2).Value = $d.B
}

# ---------------------------------------------------------------------------
# 4) Impo-ICA: revise column B values (ICA julio 24 update).
# ---------------------------------------------------------------------------
$impoIcaData = @(
    @{Row=2; B=33135.228590990002},
    @{Row=3; B=5055.0676673500002},
    @{Row=4; B=3956.5629358299998},
    @{Row=5; B=377.81753214999998},
    @{Row=6; B=720.68719937000003},
    @{Row=7; B=13278.95650133},
    @{Row=8; B=2258.2169489600001},
    @{Row=9; B=191.07411171999999},
    @{Row=10; B=663.69391570000005},
    @{Row=11; B=9488.5975255600006},
    @{Row=12; B=677.37399938999999},
    @{Row=13; B=2731.38627107},
    @{Row=14; B=563.80906234999998},
    @{Row=15; B=2167.5772087199998},
    @{Row=16; B=6968.5202222600001},
    @{Row=17; B=2753.20723116},
    @{Row=18; B=388.76236591999998},
    @{Row=19; B=3826.5506251800002},
    @{Row=20; B=3686.5149520199998},
    @{Row=21; B=278.70569719999997},
    @{Row=22; B=449.43557207999999},
    @{Row=23; B=226.47905148999999},
    @{Row=24; B=297.31865637999999},
    @{Row=25; B=978.01261603},
    @{Row=26; B=803.92100977999996},
    @{Row=27; B=652.64234906000002},
    @{Row=28; B=1271.5663743800001},
    @{Row=29; B=143.21660258}
)

$impoSheet = $wb.Worksheets.Item("Impo-ICA")
foreach ($d in $impoIcaData) {
    $impoSheet.Cells.Item($d.Row, 2).Value = $d.B
}

# ---------------------------------------------------------------------------
# 5) Tab selection moves from "Aperturas" to "BC por zonas".
# ---------------------------------------------------------------------------
$bcSheet.Activate()
